$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the active cell selection to J18
$ws.Range("J18").Select()

# Row 15: bump the base salary figure in K15, which drives B15:D15/E15 via existing formulas
$ws.Range("K15").Value = 43000

# Row 18: B18 now carries an explicit formula (was a literal 0), and C18/D18
# are rewired to chain off the prior column instead of K18 / a literal, plus K18 raised
$ws.Range("B18").Formula = "=K18*0"
$ws.Range("C18").Formula = "=B18*1.03"
$ws.Range("D18").Formula = "=C18*1.03"
$ws.Range("K18").Value = 474484

# Row 42: increase the flat figures across B:D
$ws.Range("B42:D42").Value = 3220

# Row 62: clear out the literal figures (now blank)
$ws.Range("B62:D62").ClearContents()

# Row 64: populate previously blank cells
$ws.Range("B64:D64").Value = 2000

# Row 65: populate previously blank cells
$ws.Range("B65:D65").Value = 8000
